$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change swaps the data rows for the two source files
# "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md" and
# "7d95b83c-3197-4052-a036-ceef94705299.md" across all three report sheets
# (Overview, zh-cn, de-de). Row 4 now shows the 7d95b83c record and row 5
# shows the d6eedc9f record (their "Latest ..." timestamps/filenames travel
# with the row they belong to).
# ---------------------------------------------------------------------------

# ---- Sheet "Overview" (columns A:G) ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "7d95b83c-3197-4052-a036-ceef94705299.md"
$wsOverview.Range("B4").Value = "e2e\7d95b83c-3197-4052-a036-ceef94705299.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-09-06 04:53:14"

$wsOverview.Range("A5").Value = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
$wsOverview.Range("B5").Value = "e2e\d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("G5").Value = "2016-09-06 04:52:36"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$4') {
        $hl.TextToDisplay = "e2e\7d95b83c-3197-4052-a036-ceef94705299.md"
    } elseif ($addr -eq '$B$5') {
        $hl.TextToDisplay = "e2e\d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
    }
}

# ---- Sheet "zh-cn" (columns A:P) ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "7d95b83c-3197-4052-a036-ceef94705299.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "7d95b83c-3197-4052-a036-ceef94705299.c7568932cb0bea56db9aeef1929679062ee75fff.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-06 04:53:09"

$wsZhCn.Range("A5").Value = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("G5").Value = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.1fa0a22475564896f0231cd8d45addabc8117bb8.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-06 04:52:31"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') {
        $hl.TextToDisplay = "7d95b83c-3197-4052-a036-ceef94705299.md"
    } elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
    }
}

# ---- Sheet "de-de" (columns A:P) ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "7d95b83c-3197-4052-a036-ceef94705299.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "7d95b83c-3197-4052-a036-ceef94705299.c7568932cb0bea56db9aeef1929679062ee75fff.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-06 04:53:14"

$wsDeDe.Range("A5").Value = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("G5").Value = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.1fa0a22475564896f0231cd8d45addabc8117bb8.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-06 04:52:36"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$4') {
        $hl.TextToDisplay = "7d95b83c-3197-4052-a036-ceef94705299.md"
    } elseif ($addr -eq '$A$5') {
        $hl.TextToDisplay = "d6eedc9f-8c5d-462d-b5e1-68b6d9f7ce7d.md"
    }
}

$wb.Save()
